$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trading-history rows for 2021/11/19 (serial date 44519), appended
# after the existing data (which ends at row 36).
$rows = @(
    @{ Row = 37; Stock = 3122; Action = "short"; Size = -90;  Price = 63.6 },
    @{ Row = 38; Stock = 3221; Action = "short"; Size = -61;  Price = 42.25 },
    @{ Row = 39; Stock = 6282; Action = "short"; Size = -155; Price = 35.95 },
    @{ Row = 40; Stock = 6196; Action = "short"; Size = -35;  Price = 154 },
    @{ Row = 41; Stock = 3035; Action = "short"; Size = -32;  Price = 198.5 },
    @{ Row = 42; Stock = 6411; Action = "short"; Size = -26;  Price = 272.5 },
    @{ Row = 43; Stock = 3141; Action = "short"; Size = -27;  Price = 259 },
    @{ Row = 44; Stock = 2484; Action = "long";  Size = 137;  Price = 43.65 },
    @{ Row = 45; Stock = 5351; Action = "long";  Size = 63;   Price = 93.4 },
    @{ Row = 46; Stock = 3016; Action = "long";  Size = 42;   Price = 140.5 },
    @{ Row = 47; Stock = 8069; Action = "long";  Size = 51;   Price = 115.5 },
    @{ Row = 48; Stock = 8289; Action = "long";  Size = 150;  Price = 38.35 },
    @{ Row = 49; Stock = 2340; Action = "long";  Size = 115;  Price = 51.6 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = 44519
    $dateCell.NumberFormat = $ws.Cells.Item(36, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $r.Stock
    $ws.Cells.Item($row, 3).Value = $r.Action
    $ws.Cells.Item($row, 4).Value = $r.Size
    $ws.Cells.Item($row, 5).Value = $r.Price
}

# Move the view/selection to match the state after entering the new data.
$null = $ws.Range("B37").Select()
